# Re-run of the DE density calcs: add "_basic" duplicate summary sheets
# right after the existing ones, with identical stats content.
#   area_lores      -> area_lores_basic      (copy of area_lores)
#   area_pop_sum    -> area_pop_sum_basic    (copy of area_pop_sum)

$wb = $excel.ActiveWorkbook

$wsAreaLores   = $wb.Worksheets.Item("area_lores")
$wsAreaPopSum  = $wb.Worksheets.Item("area_pop_sum")

# Duplicate "area_lores" to the end of the workbook, then rename it.
$wsAreaLores.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsAreaLoresBasic = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsAreaLoresBasic.Name = "area_lores_basic"

# Duplicate "area_pop_sum" to the end of the workbook, then rename it.
$wsAreaPopSum.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsAreaPopSumBasic = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsAreaPopSumBasic.Name = "area_pop_sum_basic"

# Keep the originally-active sheet selected, as in the source workbook.
$wsAreaLores.Activate()
